# Auto-generated edit script: updates ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# with refreshed market-price/profit values from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 237.22223
$ws.Range("J2").Value = 287.14285
$ws.Range("L2").Value = 287.14285
$ws.Range("N2").Value = -513.14285

$ws.Range("H57").Value = 51993
$ws.Range("J57").Value = 52990
$ws.Range("L57").Value = 158970
$ws.Range("N57").Value = -159968

$ws.Range("H106").Value = 5736.353
$ws.Range("J106").Value = 4339.4
$ws.Range("L106").Value = 4339.4
$ws.Range("N106").Value = -5601.4

$ws.Range("H127").Value = 512596.3
$ws.Range("I127").Value = 729137.5600000001
$ws.Range("J127").Value = 7333.3335
$ws.Range("K127").Value = 2187412.68
$ws.Range("L127").Value = 22000.0005
$ws.Range("M127").Value = -2182452.68
$ws.Range("N127").Value = -31920.0005

$ws.Range("H129").Value = 36605.95
$ws.Range("I129").Value = 50390.5
$ws.Range("K129").Value = 151171.5
$ws.Range("M129").Value = -146171.5

$ws.Range("H132").Value = 13844.167
$ws.Range("I132").Value = 12350.389
$ws.Range("J132").Value = 18325.5
$ws.Range("K132").Value = 37051.167
$ws.Range("L132").Value = 54976.5
$ws.Range("M132").Value = -34521.167
$ws.Range("N132").Value = -60036.5

$ws.Range("H135").Value = 3434.4644
$ws.Range("I135").Value = 1506.6
$ws.Range("K135").Value = 13559.4
$ws.Range("M135").Value = -11024.4

$ws.Range("H137").Value = 19495.766
$ws.Range("I137").Value = 2307.375
$ws.Range("J137").Value = 34774.332
$ws.Range("K137").Value = 6922.125
$ws.Range("L137").Value = 104322.996
$ws.Range("M137").Value = -4372.125
$ws.Range("N137").Value = -109422.996

$ws.Range("H138").Value = 3189.0293
$ws.Range("I138").Value = 9899.5
$ws.Range("J138").Value = 2769.625
$ws.Range("K138").Value = 29698.5
$ws.Range("L138").Value = 8308.875
$ws.Range("M138").Value = -24558.5
$ws.Range("N138").Value = -18588.875

$ws.Range("H141").Value = 1630.8422
$ws.Range("I141").Value = 1568.6666
$ws.Range("J141").Value = 2750
$ws.Range("K141").Value = 4705.9998
$ws.Range("L141").Value = 8250
$ws.Range("M141").Value = 474.0002000000004
$ws.Range("N141").Value = -18610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5784.0986
$ws.Range("I32").Value = 3284.7407
$ws.Range("K32").Value = 3284.7407
$ws.Range("M32").Value = -2997.7407

$ws.Range("H88").Value = 1945.6
$ws.Range("I88").Value = 1301.8334
$ws.Range("K88").Value = 1301.8334
$ws.Range("M88").Value = -895.8334

$ws.Range("H91").Value = 1945.6
$ws.Range("I91").Value = 1301.8334
$ws.Range("K91").Value = 1301.8334
$ws.Range("M91").Value = 102.1666

$ws.Range("H97").Value = 975.125
$ws.Range("I97").Value = 923.8421
$ws.Range("K97").Value = 923.8421
$ws.Range("M97").Value = -427.8421

$ws.Range("H105").Value = 225000
$ws.Range("J105").Value = 225000
$ws.Range("L105").Value = 225000
$ws.Range("N105").Value = -231988

$ws.Range("H110").Value = 10103727
$ws.Range("I110").Value = 10103727
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 10103727
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -10101682
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 7589.349
$ws.Range("I132").Value = 2206.7917
$ws.Range("J132").Value = 24813.533
$ws.Range("K132").Value = 6620.375100000001
$ws.Range("L132").Value = 74440.599
$ws.Range("M132").Value = -4090.375100000001
$ws.Range("N132").Value = -79500.599

$ws.Range("H135").Value = 165571.6
$ws.Range("J135").Value = 165571.6
$ws.Range("L135").Value = 165571.6
$ws.Range("N135").Value = -175711.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 22530.75
$ws.Range("J95").Value = 22530.75
$ws.Range("L95").Value = 22530.75
$ws.Range("N95").Value = -28022.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 48010.02
$ws.Range("I31").Value = 63167.324
$ws.Range("K31").Value = 63167.324
$ws.Range("M31").Value = -62872.324

$ws.Range("H34").Value = 48010.02
$ws.Range("I34").Value = 63167.324
$ws.Range("K34").Value = 63167.324
$ws.Range("M34").Value = -62965.324

$ws.Range("H43").Value = 17749.334
$ws.Range("J43").Value = 17749.334
$ws.Range("L43").Value = 17749.334
$ws.Range("N43").Value = -18117.334

$ws.Range("H62").Value = 3860
$ws.Range("J62").Value = 2433.3333
$ws.Range("L62").Value = 2433.3333
$ws.Range("N62").Value = -3681.3333

$ws.Range("H65").Value = 3860
$ws.Range("J65").Value = 2433.3333
$ws.Range("L65").Value = 12166.6665
$ws.Range("N65").Value = -18406.6665

$ws.Range("H101").Value = 17749.334
$ws.Range("J101").Value = 17749.334
$ws.Range("L101").Value = 17749.334
$ws.Range("N101").Value = -24239.334

$ws.Range("H134").Value = 12420.286
$ws.Range("I134").Value = 3874.25
$ws.Range("J134").Value = 15838.7
$ws.Range("K134").Value = 11622.75
$ws.Range("L134").Value = 47516.10000000001
$ws.Range("M134").Value = -9087.75
$ws.Range("N134").Value = -52586.10000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 318.7143
$ws.Range("I46").Value = 318.7143
$ws.Range("K46").Value = 956.1428999999999
$ws.Range("M46").Value = -865.1428999999999

$ws.Range("H52").Value = 1649.6
$ws.Range("I52").Value = 1999.3334
$ws.Range("K52").Value = 5998.0002
$ws.Range("M52").Value = -5732.0002

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws.Range("H113").Value = 713.61536
$ws.Range("J113").Value = 774.8889
$ws.Range("L113").Value = 2324.6667
$ws.Range("N113").Value = -6664.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 1438128.2
$ws.Range("J36").Value = 3000
$ws.Range("L36").Value = 3000
$ws.Range("N36").Value = -3970

$ws.Range("H101").Value = 22999.5
$ws.Range("J101").Value = 22999.5
$ws.Range("L101").Value = 22999.5
$ws.Range("N101").Value = -29489.5

$ws.Range("H105").Value = 111244.25
$ws.Range("J105").Value = 111244.25
$ws.Range("L105").Value = 111244.25
$ws.Range("N105").Value = -118232.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2607.111
$ws.Range("J82").Value = 3469.625
$ws.Range("L82").Value = 3469.625
$ws.Range("N82").Value = -4191.625

$ws.Range("H85").Value = 2607.111
$ws.Range("J85").Value = 3469.625
$ws.Range("L85").Value = 3469.625
$ws.Range("N85").Value = -5965.625

$ws.Range("H106").Value = 26180
$ws.Range("J106").Value = 26180
$ws.Range("L106").Value = 26180
$ws.Range("N106").Value = -28704

$ws.Range("H122").Value = 54579920
$ws.Range("I122").Value = 69437300
$ws.Range("J122").Value = 10007800
$ws.Range("K122").Value = 208311900
$ws.Range("L122").Value = 30023400
$ws.Range("M122").Value = -208309450
$ws.Range("N122").Value = -30028300

$ws.Range("H132").Value = 13839.6
$ws.Range("I132").Value = 7799.5
$ws.Range("J132").Value = 20742.572
$ws.Range("K132").Value = 23398.5
$ws.Range("L132").Value = 62227.716
$ws.Range("M132").Value = -20868.5
$ws.Range("N132").Value = -67287.716

$ws.Range("H136").Value = 45075.395
$ws.Range("I136").Value = 55507.51
$ws.Range("K136").Value = 166522.53
$ws.Range("M136").Value = -163972.53

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 51000
$ws.Range("J42").Value = 50000
$ws.Range("L42").Value = 50000
$ws.Range("N42").Value = -50756

$ws.Range("H95").Value = 32500
$ws.Range("J95").Value = 32500
$ws.Range("L95").Value = 32500
$ws.Range("N95").Value = -37992

$ws.Range("H122").Value = 720983.4399999999
$ws.Range("J122").Value = 5753.6113
$ws.Range("L122").Value = 17260.8339
$ws.Range("N122").Value = -22160.8339

$ws.Range("H132").Value = 5243.34
$ws.Range("I132").Value = 2140.9736
$ws.Range("J132").Value = 15067.5
$ws.Range("K132").Value = 6422.9208
$ws.Range("L132").Value = 45202.5
$ws.Range("M132").Value = -3892.9208
$ws.Range("N132").Value = -50262.5

Write-Host "Applied scheduled market-data update across all sheets."